$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 3.843275666666667
$ws.Range("H2").Value = 11.529827
$ws.Range("I2").Value = 0.02732998309962442
$ws.Range("J2").Value = 0.02732998309962442
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.319869
$ws.Range("N2").Value = 0.959607
$ws.Range("O2").Value = 0.17002668917608
$ws.Range("P2").Value = 0.17002668917608
$ws.Range("Q2").Value = 1.229344744221
$ws.Range("R2").Value = 11.064102697989
$ws.Range("S2").Value = 0.00464682654166736
$ws.Range("T2").Value = 0.004646826541667361

# Row 3
$ws.Range("G3").Value = 3.843275666666667
$ws.Range("H3").Value = 11.529827
$ws.Range("I3").Value = 0.02732998309962442
$ws.Range("J3").Value = 0.02732998309962442
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.3653903333333333
$ws.Range("N3").Value = 1.096171
$ws.Range("O3").Value = 0.1942235997661884
$ws.Range("P3").Value = 0.1942235997661884
$ws.Range("Q3").Value = 1.404295776935222
$ws.Range("R3").Value = 12.638661992417
$ws.Range("S3").Value = 0.005308127699158146
$ws.Range("T3").Value = 0.005308127699158147

# Row 4
$ws.Range("G4").Value = 3.843275666666667
$ws.Range("H4").Value = 11.529827
$ws.Range("I4").Value = 0.02732998309962442
$ws.Range("J4").Value = 0.02732998309962442
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 1.196027666666667
$ws.Range("N4").Value = 3.588083
$ws.Range("O4").Value = 0.6357497110577316
$ws.Range("P4").Value = 0.6357497110577316
$ws.Range("Q4").Value = 4.596664027960111
$ws.Range("R4").Value = 41.369976251641
$ws.Range("S4").Value = 0.01737502885879891
$ws.Range("T4").Value = 0.01737502885879891

# Row 5
$ws.Range("G5").Value = 9.653191
$ws.Range("H5").Value = 28.959573
$ws.Range("I5").Value = 0.06864497105310771
$ws.Range("J5").Value = 0.06864497105310771
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.319869
$ws.Range("N5").Value = 0.959607
$ws.Range("O5").Value = 0.17002668917608
$ws.Range("P5").Value = 0.17002668917608
$ws.Range("Q5").Value = 3.087756551979
$ws.Range("R5").Value = 27.789808967811
$ws.Range("S5").Value = 0.01167147715674775
$ws.Range("T5").Value = 0.01167147715674775

# Row 6
$ws.Range("G6").Value = 9.653191
$ws.Range("H6").Value = 28.959573
$ws.Range("I6").Value = 0.06864497105310771
$ws.Range("J6").Value = 0.06864497105310771
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.3653903333333333
$ws.Range("N6").Value = 1.096171
$ws.Range("O6").Value = 0.1942235997661884
$ws.Range("P6").Value = 0.1942235997661884
$ws.Range("Q6").Value = 3.527182677220333
$ws.Range("R6").Value = 31.744644094983
$ws.Range("S6").Value = 0.01333247338378038
$ws.Range("T6").Value = 0.01333247338378038

# Row 7
$ws.Range("G7").Value = 9.653191
$ws.Range("H7").Value = 28.959573
$ws.Range("I7").Value = 0.06864497105310771
$ws.Range("J7").Value = 0.06864497105310771
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 1.196027666666667
$ws.Range("N7").Value = 3.588083
$ws.Range("O7").Value = 0.6357497110577316
$ws.Range("P7").Value = 0.6357497110577316
$ws.Range("Q7").Value = 11.54548350761767
$ws.Range("R7").Value = 103.909351568559
$ws.Range("S7").Value = 0.04364102051257958
$ws.Range("T7").Value = 0.04364102051257958

# Row 8
$ws.Range("G8").Value = 0.5927003333333334
$ws.Range("H8").Value = 1.778101
$ws.Range("I8").Value = 0.004214761442597993
$ws.Range("J8").Value = 0.004214761442597993
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.319869
$ws.Range("N8").Value = 0.959607
$ws.Range("O8").Value = 0.17002668917608
$ws.Range("P8").Value = 0.17002668917608
$ws.Range("Q8").Value = 0.189586462923
$ws.Range("R8").Value = 1.706278166307
$ws.Range("S8").Value = 0.0007166219337519354
$ws.Range("T8").Value = 0.0007166219337519354

# Row 9
$ws.Range("G9").Value = 0.5927003333333334
$ws.Range("H9").Value = 1.778101
$ws.Range("I9").Value = 0.004214761442597993
$ws.Range("J9").Value = 0.004214761442597993
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.3653903333333333
$ws.Range("N9").Value = 1.096171
$ws.Range("O9").Value = 0.1942235997661884
$ws.Range("P9").Value = 0.1942235997661884
$ws.Range("Q9").Value = 0.2165669723634445
$ws.Range("R9").Value = 1.949102751271
$ws.Range("S9").Value = 0.0008186061395371155
$ws.Range("T9").Value = 0.0008186061395371155

# Row 10
$ws.Range("G10").Value = 0.5927003333333334
$ws.Range("H10").Value = 1.778101
$ws.Range("I10").Value = 0.004214761442597993
$ws.Range("J10").Value = 0.004214761442597993
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 1.196027666666667
$ws.Range("N10").Value = 3.588083
$ws.Range("O10").Value = 0.6357497110577316
$ws.Range("P10").Value = 0.6357497110577316
$ws.Range("Q10").Value = 0.7088859967092223
$ws.Range("R10").Value = 6.379973970383001
$ws.Range("S10").Value = 0.002679533369308942
$ws.Range("T10").Value = 0.002679533369308942

# Row 11
$ws.Range("G11").Value = 126.5357156666666
$ws.Range("H11").Value = 379.6071469999999
$ws.Range("I11").Value = 0.8998102844046698
$ws.Range("J11").Value = 0.8998102844046699
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.319869
$ws.Range("N11").Value = 0.959607
$ws.Range("O11").Value = 0.17002668917608
$ws.Range("P11").Value = 0.17002668917608
$ws.Range("Q11").Value = 40.47485283458099
$ws.Range("R11").Value = 364.2736755112289
$ws.Range("S11").Value = 0.1529917635439129
$ws.Range("T11").Value = 0.1529917635439129

# Row 12
$ws.Range("G12").Value = 126.5357156666666
$ws.Range("H12").Value = 379.6071469999999
$ws.Range("I12").Value = 0.8998102844046698
$ws.Range("J12").Value = 0.8998102844046699
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.3653903333333333
$ws.Range("N12").Value = 1.096171
$ws.Range("O12").Value = 0.1942235997661884
$ws.Range("P12").Value = 0.1942235997661884
$ws.Range("Q12").Value = 46.23492732601521
$ws.Range("R12").Value = 416.114345934137
$ws.Range("S12").Value = 0.1747643925437128
$ws.Range("T12").Value = 0.1747643925437128

# Row 13
$ws.Range("G13").Value = 126.5357156666666
$ws.Range("H13").Value = 379.6071469999999
$ws.Range("I13").Value = 0.8998102844046698
$ws.Range("J13").Value = 0.8998102844046699
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 1.196027666666667
$ws.Range("N13").Value = 3.588083
$ws.Range("O13").Value = 0.6357497110577316
$ws.Range("P13").Value = 0.6357497110577316
$ws.Range("Q13").Value = 151.3402167588001
$ws.Range("R13").Value = 1362.061950829201
$ws.Range("S13").Value = 0.5720541283170442
$ws.Range("T13").Value = 0.5720541283170442
